# Update task data used in testing (row 5 of Sheet1) and move the
# active selection to D5, matching the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 3
$ws.Range("F5").Value = 3
$ws.Range("H5").Value = 46

$ws.Range("D5").Select()
